$d = $word.ActiveDocument

# Locate the anchor paragraph ("LOB1039: Física Experimental III (Requisito fraco)")
# and remove the three paragraphs that directly follow it:
#   1) an empty "Normal" paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. ..."
# while leaving the subsequent empty paragraph and the page-break paragraph intact.

$anchorText = "LOB1039: Física Experimental III (Requisito fraco)"

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $anchorText) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    $firstToRemove = $d.Paragraphs.Item($anchorIndex + 1)
    $lastToRemove = $d.Paragraphs.Item($anchorIndex + 3)

    $start = $firstToRemove.Range.Start
    $end = $lastToRemove.Range.End

    $killRange = $d.Range($start, $end)
    $killRange.Delete()
}
